{"js": "const body = context.document.body;\n\n// Heading3 paragraph: \"[03/07/24]\" with bookmark \"section-2\"\nconst heading = body.insertParagraph(\"[03/07/24]\", Word.InsertLocation.end);\nheading.style = \"Heading 3\";\nheading.getRange().insertBookmark(\"section-2\");\n\n// FirstParagraph paragraph\nconst firstPara = body.insertParagraph(\n  \"Denna dag arbeta vi p\u00e5 att experimentera med v\u00e5rt torn, vi f\u00f6rs\u00f6kte g\u00f6ra det l\u00e4ngre med ben, och starkare med tejp, tyv\u00e4rr gick b\u00e5de experimenten r\u00e4tt d\u00e5ligt.\",\n  Word.InsertLocation.end\n);\nfirstPara.style = \"First Paragraph\";\n\n// BodyText paragraph\nconst bodyPara = body.insertParagraph(\n  \"Till slut f\u00f6rst\u00f6rdes v\u00e5rt torns stabilitet p.g.a. vi hade f\u00f6r mycket vikt p\u00e5.\",\n  Word.InsertLocation.end\n);\nbodyPara.style = \"Body Text\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Move to the very end of the document content.\n$r = $d.Content\n$r.Collapse(0)\n\n# --- Heading3 paragraph: \"[03/07/24]\" with bookmark \"section-2\" ---\n$r.InsertParagraphAfter()\n$r.Collapse(0)\n$p1 = $d.Paragraphs.Last\n$p1.Style = \"Heading 3\"\n$p1.Range.Text = \"[03/07/24]\"\n$d.Bookmarks.Add(\"section-2\", $p1.Range)\n\n# --- FirstParagraph paragraph ---\n$r = $d.Content\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n$r.Collapse(0)\n$p2 = $d.Paragraphs.Last\n$p2.Style = \"First Paragraph\"\n$p2.Range.Text = \"Denna dag arbeta vi p\u00e5 att experimentera med v\u00e5rt torn, vi f\u00f6rs\u00f6kte g\u00f6ra det l\u00e4ngre med ben, och starkare med tejp, tyv\u00e4rr gick b\u00e5de experimenten r\u00e4tt d\u00e5ligt.\"\n\n# --- BodyText paragraph ---\n$r = $d.Content\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n$r.Collapse(0)\n$p3 = $d.Paragraphs.Last\n$p3.Style = \"Body Text\"\n$p3.Range.Text = \"Till slut f\u00f6rst\u00f6rdes v\u00e5rt torns stabilitet p.g.a. vi hade f\u00f6r mycket vikt p\u00e5.\"\n"}
